$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.461.85'
$ws.Range("E2").Value = '  -2.96%  '
$ws.Range("D3").Value = '3.315.48'
$ws.Range("E3").Value = '  -4.95%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '548.42'
$ws.Range("E5").Value = '  -0.83%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.84'
$ws.Range("E6").Value = '  -3.55%  '
$ws.Range("E7").Value = '  -4.68%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("D9").Value = '3.305.77'
$ws.Range("E9").Value = '  -5.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.611'
$ws.Range("E10").Value = '  -3.35%  '
$ws.Range("E11").Value = '  -2.02%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.10'
$ws.Range("E12").Value = '  -1.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000265'
$ws.Range("E13").Value = '  -2.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '8.89'
$ws.Range("E14").Value = '  -2.88%  '
$ws.Range("D15").Value = '3.841.14'
$ws.Range("E15").Value = '  -5.00%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.18'
$ws.Range("E16").Value = '  -1.11%  '
$ws.Range("E17").Value = '  -3.35%  '
$ws.Range("D18").Value = '3.310.97'
$ws.Range("E18").Value = '  -5.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.68'
$ws.Range("E19").Value = '  -4.16%  '
$ws.Range("D20").Value = '63.393.71'
$ws.Range("E20").Value = '  -3.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.968'
$ws.Range("E21").Value = '  -2.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '424.06'
$ws.Range("E22").Value = '  +2.60%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.43'
$ws.Range("E23").Value = '  +8.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.05'
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("B25").Value = 'InternetComputer(DFINITY)'
$ws.Range("C25").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '13.30'
$ws.Range("E25").Value = '  +4.39%  '
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '83.07'
$ws.Range("E26").Value = '  -3.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.60'
$ws.Range("E27").Value = '  -1.64%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.72'
$ws.Range("E28").Value = '  -4.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.66'
$ws.Range("E29").Value = '  -3.92%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '29.16'
$ws.Range("E30").Value = '  -3.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.42'
$ws.Range("E31").Value = '  -0.60%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.35'
$ws.Range("E32").Value = '  -2.56%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '575.94'
$ws.Range("E33").Value = '  -6.07%  '
$ws.Range("E34").Value = '  -3.37%  '
$ws.Range("E35").Value = '  -2.20%  '
$ws.Range("E36").Value = '  -0.10%  '
$ws.Range("E37").Value = '  -1.61%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.48'
$ws.Range("E38").Value = '  +6.97%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.10'
$ws.Range("E39").Value = '  -5.32%  '
$ws.Range("D40").Value = '0.0₃0738'
$ws.Range("E40").Value = '  -6.54%  '
$ws.Range("E41").Value = '  -4.33%  '
$ws.Range("D42").Value = '3.122.66'
$ws.Range("E42").Value = '  -7.29%  '
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.78'
$ws.Range("E44").Value = '  -2.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.19'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0402'
$ws.Range("E46").Value = '  -3.09%  '
$ws.Range("E47").Value = '  -3.35%  '
$ws.Range("E48").Value = '  -6.32%  '
$ws.Range("E49").Value = '  -3.90%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.27'
$ws.Range("E50").Value = '  -2.98%  '
$ws.Range("E51").Value = '  -4.67%  '
